# Apply updated cryptocurrency price/volume data (and the row 32/33 coin swap)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $value) {
    # Force the cell to stay a text value (matches the source sheet, which
    # stores every Coin/Link/Price/Volume cell as a string) even when the
    # text looks like a plain number (e.g. "290.71"), then restore the
    # default "Normal" style so we do not leave a stray number-format behind.
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell 2 4 '42.003.66'
Set-TextCell 2 5 '  -9.17%  '

# Row 3
Set-TextCell 3 4 '2.473.04'
Set-TextCell 3 5 '  -5.00%  '

# Row 4
Set-TextCell 4 4 '0.999'
Set-TextCell 4 5 '  -0.15%  '

# Row 5
Set-TextCell 5 4 '290.71'
Set-TextCell 5 5 '  -5.24%  '

# Row 6
Set-TextCell 6 4 '90.99'
Set-TextCell 6 5 '  -8.37%  '

# Row 7
Set-TextCell 7 4 '0.561'
Set-TextCell 7 5 '  -6.50%  '

# Row 8
Set-TextCell 8 4 '1.00'
Set-TextCell 8 5 '  +0.14%  '

# Row 9
Set-TextCell 9 4 '0.534'
Set-TextCell 9 5 '  -7.37%  '

# Row 10
Set-TextCell 10 4 '35.33'
Set-TextCell 10 5 '  -10.19%  '

# Row 11
Set-TextCell 11 4 '0.0784'
Set-TextCell 11 5 '  -6.84%  '

# Row 12
Set-TextCell 12 4 '7.51'
Set-TextCell 12 5 '  -7.16%  '

# Row 13
Set-TextCell 13 4 '0.106'
Set-TextCell 13 5 '  +0.22%  '

# Row 14
Set-TextCell 14 4 '2.871.53'
Set-TextCell 14 5 '  -4.52%  '

# Row 15
Set-TextCell 15 4 '2.504.73'
Set-TextCell 15 5 '  -3.97%  '

# Row 16
Set-TextCell 16 4 '0.846'
Set-TextCell 16 5 '  -7.69%  '

# Row 17
Set-TextCell 17 4 '13.74'
Set-TextCell 17 5 '  -7.70%  '

# Row 18
Set-TextCell 18 4 '42.023.77'
Set-TextCell 18 5 '  -9.48%  '

# Row 19
Set-TextCell 19 4 '0.0₃0940'
Set-TextCell 19 5 '  -6.67%  '

# Row 20
Set-TextCell 20 4 '6.38'
Set-TextCell 20 5 '  -4.62%  '

# Row 21
Set-TextCell 21 4 '11.96'
Set-TextCell 21 5 '  -7.10%  '

# Row 22
Set-TextCell 22 4 '71.11'
Set-TextCell 22 5 '  -0.19%  '

# Row 23
Set-TextCell 23 4 '253.01'
Set-TextCell 23 5 '  -6.94%  '

# Row 24
Set-TextCell 24 4 '2.82'
Set-TextCell 24 5 '  -6.81%  '

# Row 25
Set-TextCell 25 4 '2.06'
Set-TextCell 25 5 '  -4.40%  '

# Row 26
Set-TextCell 26 4 '27.83'
Set-TextCell 26 5 '  -3.97%  '

# Row 27
Set-TextCell 27 4 '0.998'
Set-TextCell 27 5 '  -0.22%  '

# Row 28
Set-TextCell 28 4 '2.20'
Set-TextCell 28 5 '  +0.03%  '

# Row 29
Set-TextCell 29 4 '9.73'
Set-TextCell 29 5 '  -7.75%  '

# Row 30
Set-TextCell 30 4 '35.85'
Set-TextCell 30 5 '  -6.74%  '

# Row 31
Set-TextCell 31 4 '5.85'
Set-TextCell 31 5 '  -7.21%  '

# Row 32
Set-TextCell 32 2 'Monero'
Set-TextCell 32 3 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell 32 4 '150.11'
Set-TextCell 32 5 '  -0.59%  '

# Row 33
Set-TextCell 33 2 'LidoDAOToken'
Set-TextCell 33 3 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextCell 33 4 '3.39'
Set-TextCell 33 5 '  -6.74%  '

# Row 34
Set-TextCell 34 4 '2.12'
Set-TextCell 34 5 '  -4.78%  '

# Row 35
Set-TextCell 35 4 '2.69'
Set-TextCell 35 5 '  -5.71%  '

# Row 36
Set-TextCell 36 4 '0.0777'
Set-TextCell 36 5 '  -6.82%  '

# Row 37
Set-TextCell 37 4 '0.111'
Set-TextCell 37 5 '  -8.95%  '

# Row 38
Set-TextCell 38 5 '  -4.60%  '

# Row 39
Set-TextCell 39 4 '23.08'
Set-TextCell 39 5 '  -0.59%  '

# Row 40
Set-TextCell 40 4 '16.23'
Set-TextCell 40 5 '  +2.62%  '

# Row 41
Set-TextCell 41 4 '3.34'
Set-TextCell 41 5 '  -7.39%  '

# Row 42
Set-TextCell 42 4 '0.0300'
Set-TextCell 42 5 '  -8.76%  '

# Row 43
Set-TextCell 43 4 '3.71'
Set-TextCell 43 5 '  -8.37%  '

# Row 44
Set-TextCell 44 4 '1.984.11'
Set-TextCell 44 5 '  -6.12%  '

# Row 45
Set-TextCell 45 5 '  +0.07%  '

# Row 46
Set-TextCell 46 4 '1.61'
Set-TextCell 46 5 '  +4.54%  '

# Row 47
Set-TextCell 47 4 '83.71'
Set-TextCell 47 5 '  -10.01%  '

# Row 48
Set-TextCell 48 4 '8.76'
Set-TextCell 48 5 '  -8.26%  '

# Row 49
Set-TextCell 49 4 '2.752.24'
Set-TextCell 49 5 '  -3.91%  '

# Row 50
Set-TextCell 50 4 '100.64'
Set-TextCell 50 5 '  -7.11%  '

# Row 51
Set-TextCell 51 4 '0.183'
Set-TextCell 51 5 '  -8.50%  '
